# Techdata1.xlsx edit
# "Excel extraction is working, simple selection of the first sheet."
#
# Adds a ProductionDate/ProductionDateText pair of rows to "First sheet",
# reformats the WorkWidth value to a thousands-separated number, and leaves
# the cursor/selection parked on the new date cell - mirroring a short,
# real authoring session against the workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("First sheet")

# --- Row 8: ProductionDate (date-formatted value) ---------------------
$ws1.Range("A8").Value = "ProductionDate"
$ws1.Range("B8").NumberFormat = "mm-dd-yy"
$ws1.Range("B8").Value = 42429

# --- Row 7: ProductionDateText (plain text date look-alike) -----------
$ws1.Range("A7").Value = "ProductionDateText"
$ws1.Range("B7").NumberFormat = "@"
$ws1.Range("B7").Value = "2016-02-01"

# --- WorkWidth value bumped up and shown with a thousands separator ---
$ws1.Range("B2").NumberFormat = "#,##0"
$ws1.Range("B2").Value = 6000

# --- Make sure the first sheet is the active/selected tab, cursor on
#     the newly entered date cell (matches the final sheetView selection)
$ws1.Select() | Out-Null
$ws1.Range("B8").Select() | Out-Null
